# "edit readme and todo" - update the CodePen example link used in the
# Options table (column D, "Example") for the rows that previously had an
# empty CodePen link, and move the sheet's active selection to E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newLink = "[CodePen](https://codepen.io/maiCoding/pen/OajRdb)"

$ws.Range("D2").Value = $newLink
$ws.Range("D7").Value = $newLink
$ws.Range("D8").Value = $newLink
$ws.Range("D9").Value = $newLink

$ws.Activate()
$ws.Range("E10").Select()
